$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.30%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.75%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.591"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.71%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08166"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.18%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.775"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.42%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.48%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.904"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-6.55%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9456"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.68%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1225"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.30%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1934"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.44%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09828"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.79%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04485"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "14.51%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1069"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.82%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001275"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.97%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006030"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.62%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.502"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.69%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.726"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.88%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.95%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04395"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.54%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001240"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.62%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004345"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.81%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.01%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004007"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "31.59%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02813"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05726"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.42%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007905"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.38%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009808"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "9.78%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.46%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002101"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.80%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009790"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-16.69%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007319"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.39%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.49%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003432"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "7.97%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002280"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.20%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.49%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002009"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.49%"
